$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 27 corresponds to @nitsurathkla
# current_phase: 1 -> 2
$ws.Range("D27").Value = 2

# last_action_date: (empty) -> timestamp
$ws.Range("E27").Value = "2026-02-20T05:57:45.696563+00:00"

# reactions_count: 0 -> 1
$ws.Range("H27").Value = 1

# replies_count: 0 -> 1
$ws.Range("I27").Value = 1

# reacted_message_ids: [] -> [7225]
$ws.Range("L27").Value = "[7225]"

# replied_message_ids: [] -> [7224]
$ws.Range("M27").Value = "[7224]"
